$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 301
$ws.Range("I19").Value = 300
$ws.Range("J19").Value = 302
$ws.Range("K19").Value = 300
$ws.Range("L19").Value = 302
$ws.Range("M19").Value = -125
$ws.Range("N19").Value = -652
$ws.Range("H64").Value = 4912.5
$ws.Range("J64").Value = 3150
$ws.Range("L64").Value = 3150
$ws.Range("N64").Value = -3646
$ws.Range("H67").Value = 4912.5
$ws.Range("J67").Value = 3150
$ws.Range("L67").Value = 3150
$ws.Range("N67").Value = -4866
$ws.Range("H131").Value = 5841.1816
$ws.Range("I131").Value = 1950.6
$ws.Range("J131").Value = 9083.333000000001
$ws.Range("K131").Value = 5851.799999999999
$ws.Range("L131").Value = 27249.999
$ws.Range("M131").Value = -811.7999999999993
$ws.Range("N131").Value = -37329.999
$ws.Range("H135").Value = 2577.5
$ws.Range("J135").Value = 3204.4
$ws.Range("L135").Value = 28839.6
$ws.Range("N135").Value = -33909.60000000001
$ws.Range("H138").Value = 4048.204
$ws.Range("I138").Value = 1798.8889
$ws.Range("J138").Value = 4554.3
$ws.Range("K138").Value = 5396.6667
$ws.Range("L138").Value = 13662.9
$ws.Range("M138").Value = -256.6666999999998
$ws.Range("N138").Value = -23942.9

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4997.6
$ws.Range("I61").Value = 4997.5
$ws.Range("K61").Value = 4997.5
$ws.Range("M61").Value = -4785.5
$ws.Range("H95").Value = 48971.168
$ws.Range("J95").Value = 48971.168
$ws.Range("L95").Value = 48971.168
$ws.Range("N95").Value = -54463.168
$ws.Range("H132").Value = 1695.2333
$ws.Range("I132").Value = 1365.1111
$ws.Range("K132").Value = 4095.3333
$ws.Range("M132").Value = -1565.3333
$ws.Range("H136").Value = 4997.6
$ws.Range("I136").Value = 4997.5
$ws.Range("K136").Value = 14992.5
$ws.Range("M136").Value = -12442.5

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2788
$ws.Range("I20").Value = 726.25
$ws.Range("K20").Value = 726.25
$ws.Range("M20").Value = -479.25
$ws.Range("H86").Value = 1563.7
$ws.Range("I86").Value = 1841
$ws.Range("J86").Value = 1147.75
$ws.Range("K86").Value = 1841
$ws.Range("L86").Value = 1147.75
$ws.Range("M86").Value = -718
$ws.Range("N86").Value = -3393.75
$ws.Range("H89").Value = 1563.7
$ws.Range("I89").Value = 1841
$ws.Range("J89").Value = 1147.75
$ws.Range("K89").Value = 9205
$ws.Range("L89").Value = 5738.75
$ws.Range("M89").Value = -3589
$ws.Range("N89").Value = -16970.75
$ws.Range("H107").Value = 635.4231
$ws.Range("I107").Value = 627.04346
$ws.Range("K107").Value = 627.04346
$ws.Range("M107").Value = 1292.95654
$ws.Range("H134").Value = 2464.0715
$ws.Range("I134").Value = 2011.7
$ws.Range("K134").Value = 6035.1
$ws.Range("M134").Value = -3500.1

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 5822.5713
$ws.Range("I58").Value = 4812.7144
$ws.Range("J58").Value = 6832.4287
$ws.Range("K58").Value = 4812.7144
$ws.Range("L58").Value = 6832.4287
$ws.Range("M58").Value = -4609.7144
$ws.Range("N58").Value = -7238.4287
$ws.Range("H134").Value = 2870.258
$ws.Range("I134").Value = 2070.1
$ws.Range("J134").Value = 4325.091
$ws.Range("K134").Value = 6210.299999999999
$ws.Range("L134").Value = 12975.273
$ws.Range("M134").Value = -3675.299999999999
$ws.Range("N134").Value = -18045.273
$ws.Range("H136").Value = 5822.5713
$ws.Range("I136").Value = 4812.7144
$ws.Range("J136").Value = 6832.4287
$ws.Range("K136").Value = 14438.1432
$ws.Range("L136").Value = 20497.2861
$ws.Range("M136").Value = -11888.1432
$ws.Range("N136").Value = -25597.2861

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 62557.688
$ws.Range("J2").Value = 95
$ws.Range("L2").Value = 570
$ws.Range("N2").Value = -796
$ws.Range("H34").Value = 1160.0769
$ws.Range("I34").Value = 654.125
$ws.Range("K34").Value = 1962.375
$ws.Range("M34").Value = -1878.375
$ws.Range("H47").Value = 95
$ws.Range("J47").Value = 90
$ws.Range("L47").Value = 270
$ws.Range("N47").Value = -1132
$ws.Range("H129").Value = 2686
$ws.Range("I129").Value = 1997.75
$ws.Range("J129").Value = 2991.889
$ws.Range("K129").Value = 5993.25
$ws.Range("L129").Value = 8975.667000000001
$ws.Range("M129").Value = -993.25
$ws.Range("N129").Value = -18975.667
$ws.Range("H137").Value = 4378.579
$ws.Range("I137").Value = 2959.2
$ws.Range("J137").Value = 5955.6665
$ws.Range("K137").Value = 8877.599999999999
$ws.Range("L137").Value = 17866.9995
$ws.Range("M137").Value = -3777.599999999999
$ws.Range("N137").Value = -28066.9995

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H109").Value = 42000
$ws.Range("J109").Value = 42000
$ws.Range("L109").Value = 42000
$ws.Range("N109").Value = -44080
$ws.Range("H132").Value = 1879.625
$ws.Range("J132").Value = 4496
$ws.Range("L132").Value = 13488
$ws.Range("N132").Value = -18548

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 1590.3334
$ws.Range("I40").Value = 1590.3334
$ws.Range("K40").Value = 1590.3334
$ws.Range("M40").Value = -1454.3334
$ws.Range("H55").Value = 345.75
$ws.Range("I55").Value = 359.46155
$ws.Range("K55").Value = 359.46155
$ws.Range("M55").Value = -186.46155
$ws.Range("H122").Value = 4021.3635
$ws.Range("I122").Value = 4048.5
$ws.Range("J122").Value = 3750
$ws.Range("K122").Value = 12145.5
$ws.Range("L122").Value = 11250
$ws.Range("M122").Value = -9695.5
$ws.Range("N122").Value = -16150
$ws.Range("H132").Value = 4835.3335
$ws.Range("I132").Value = 5003
$ws.Range("K132").Value = 15009
$ws.Range("M132").Value = -12479

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 4071.75
$ws.Range("I122").Value = 4795.8335
$ws.Range("J122").Value = 1899.5
$ws.Range("K122").Value = 14387.5005
$ws.Range("L122").Value = 5698.5
$ws.Range("M122").Value = -11937.5005
$ws.Range("N122").Value = -10598.5
$ws.Range("H125").Value = 127300
$ws.Range("J125").Value = 127300
$ws.Range("L125").Value = 127300
$ws.Range("N125").Value = -137140
$ws.Range("H126").Value = 1703.1052
$ws.Range("I126").Value = 1096.6
$ws.Range("J126").Value = 3977.5
$ws.Range("K126").Value = 3289.8
$ws.Range("L126").Value = 11932.5
$ws.Range("M126").Value = -819.7999999999997
$ws.Range("N126").Value = -16872.5
$ws.Range("H132").Value = 1844.25
$ws.Range("I132").Value = 1536.2858
$ws.Range("K132").Value = 4608.857400000001
$ws.Range("M132").Value = -2078.857400000001
